# The document contains several paragraphs whose visible text is the
# literal string "<id>p142v_N</id>" (an XML-ish tag wrapped around a
# page/id value). Each such paragraph was left over-split into three
# runs -- one for the "<id>" open tag, one (differently formatted) for
# the inner value, and one for the "</id>" close tag -- even though the
# open/close tag runs share identical run formatting. This script
# merges every such split back into a single run per paragraph, i.e.
# "<id>" + "p142v_N" + "</id>" becomes one run "<id>p142v_N</id>" that
# keeps the "<id>"/"</id>" run's formatting (Courier New monospace,
# brown color, 9pt), matching how the rest of the document's tag runs
# ("<div>", "<head>", ...) are already stored as single runs.

$d = $word.ActiveDocument

$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs($i)
    $pRange = $p.Range
    $fullText = $pRange.Text
    # Strip the trailing paragraph-mark character(s) to get the visible text.
    $visible = $fullText.TrimEnd([char]13, [char]7)

    if (-not ($visible.StartsWith("<id>") -and $visible.EndsWith("</id>"))) {
        continue
    }
    # Need room for an inner value between the open/close tags.
    if ($visible.Length -le ("<id>".Length + "</id>".Length)) {
        continue
    }

    $paraStart = $pRange.Start
    $paraEndNoMark = $pRange.End - 1

    # Locate the end of the leading "<id>" run inside this paragraph.
    $openRange = $d.Range($paraStart, $paraEndNoMark)
    $openFound = $openRange.Find.Execute("<id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $openFound) {
        continue
    }
    $openEnd = $openRange.End

    # Already a single run covering the whole tag -- nothing to merge.
    if ($openEnd -ge $paraEndNoMark) {
        continue
    }

    # Range for everything after "<id>" up to (not including) the
    # paragraph mark -- this is the inner value run plus the "</id>" run.
    $tailRange = $d.Range($openEnd, $paraEndNoMark)
    $tailText = $tailRange.Text

    # Range for the leading "<id>" run; inserting after it picks up its
    # run formatting (Courier New / brown / 9pt) for the merged text.
    $leadRange = $d.Range($paraStart, $openEnd)

    $tailRange.Delete()
    $leadRange.InsertAfter($tailText)
}
